# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for the "展览" sheet (all rows changed).
$zhanlanUpdates = @{
    3  = 810
    4  = 272
    8  = 2039
    9  = 7588
    10 = 903
    11 = 417
    12 = 349
    13 = 133
    16 = 7737
    17 = 305
    18 = 1341
    19 = 149
    22 = 146
    23 = 307
    26 = 18
    29 = 410
    30 = 612
    32 = 91
    35 = 37
}

# Row -> new F-column value for the "全部类型" sheet (same as above,
# except row 8, whose F value was already up to date on this sheet).
$quanbuUpdates = @{
    3  = 810
    4  = 272
    9  = 7588
    10 = 903
    11 = 417
    12 = 349
    13 = 133
    16 = 7737
    17 = 305
    18 = 1341
    19 = 149
    22 = 146
    23 = 307
    26 = 18
    29 = 410
    30 = 612
    32 = 91
    35 = 37
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
